{"js": "// Add \"Num\u00e9ro agr\u00e9ment : {{ etablissement.numero_agrement }}\" as a new\n// paragraph right after the \"N\u00b0 SIRET : {{ etablissement.siret }}\"\n// paragraph in the TIAC \"etablissement\" block, matching the existing\n// Corpsdetexte styling used by the sibling fields.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet siretParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"N\u00b0 SIRET\") !== -1) {\n    siretParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!siretParagraph) {\n  throw new Error('Paragraph containing \"N\u00b0 SIRET\" was not found.');\n}\n\n// Build the raw paragraph markup exactly as authored (three runs, because\n// the middle run only carries the space and keeps a distinct/empty complex\n// -script font reference, the same shape Word leaves behind after manual\n// editing of the placeholder text).\nconst newParagraphXml =\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:pStyle w:val=\"Corpsdetexte\"/>' +\n      '<w:rPr>' +\n        '<w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/>' +\n        '<w:sz w:val=\"20\"/>' +\n        '<w:szCs w:val=\"20\"/>' +\n      '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:cs=\"Calibri\" w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/>' +\n        '<w:sz w:val=\"20\"/>' +\n        '<w:szCs w:val=\"20\"/>' +\n      '</w:rPr>' +\n      '<w:t>Num\\u00e9ro agr\\u00e9ment : {{ etablissement.numero_agrement</w:t>' +\n    '</w:r>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:cs=\"\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/>' +\n        '<w:sz w:val=\"20\"/>' +\n        '<w:szCs w:val=\"20\"/>' +\n      '</w:rPr>' +\n      '<w:t xml:space=\"preserve\"> </w:t>' +\n    '</w:r>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:cs=\"Calibri\" w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/>' +\n        '<w:sz w:val=\"20\"/>' +\n        '<w:szCs w:val=\"20\"/>' +\n      '</w:rPr>' +\n      '<w:t>}}</w:t>' +\n    '</w:r>' +\n  '</w:p>';\n\nconst flatOpcXml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' + newParagraphXml + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\n// Collapse to a caret right after the SIRET paragraph (including its\n// paragraph mark) and insert the new paragraph there, so it lands between\n// the SIRET line and the \"Autre identifiant\" line.\nconst insertionPoint = siretParagraph.getRange(\"End\");\ninsertionPoint.insertOoxml(flatOpcXml, \"After\");\nawait context.sync();\n", "ps1": "# Add \"Num\u00e9ro agr\u00e9ment : {{ etablissement.numero_agrement }}\" as a new\n# paragraph right after the \"N\u00b0 SIRET : {{ etablissement.siret }}\"\n# paragraph in the TIAC \"etablissement\" block, matching the existing\n# Corpsdetexte styling used by the sibling fields.\n\n$d = $word.ActiveDocument\n\n$siretParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*SIRET*\") {\n        $siretParagraph = $p\n        break\n    }\n}\nif ($null -eq $siretParagraph) {\n    throw 'Paragraph containing \"N\u00b0 SIRET\" was not found.'\n}\n\n# Open up a fresh (empty) paragraph right after the SIRET line, inheriting\n# its paragraph/run formatting (Corpsdetexte style, Calibri 10pt).\n$siretParagraph.Range.InsertParagraphAfter()\n$newParagraph = $siretParagraph.Next()\n\n# Fill that new paragraph with the exact markup (three runs, because the\n# middle run only carries the space and keeps a distinct/empty complex\n# -script font reference, the same shape Word leaves behind after manual\n# editing of the placeholder text).\n$newParagraphXml = @'\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"Corpsdetexte\"/><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Calibri\" w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t>Num\u00e9ro agr\u00e9ment : {{ etablissement.numero_agrement</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Calibri\" w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n$newParagraph.Range.InsertXML($newParagraphXml)\n"}
